$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# D-column values that parse as plain decimals get a leading apostrophe so
# Excel keeps them as text instead of coercing them to a Number; the Style
# is then reset to Normal so no extra number-format/quote-prefix styling is
# left attached to the cell (matches the source cells, which carry no style).

$ws.Range("D2").Value = '61.629.76'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '2.924.82'
$ws.Range("E3").Value = '  +0.61%  '
$c = $ws.Range("D4")
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$c.Value = '''595.85'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.26%  '
$c = $ws.Range("D6")
$c.Value = '''141.44'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("D8").Value = '2.921.92'
$ws.Range("E8").Value = '  +0.57%  '
$c = $ws.Range("D9")
$c.Value = '''0.498'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.30%  '
$c = $ws.Range("D10")
$c.Value = '''7.17'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.14%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("E13").Value = '  -0.91%  '
$c = $ws.Range("D14")
$c.Value = '''33.03'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").Value = '3.409.07'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '61.471.13'
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = '2.903.89'
$ws.Range("E19").Value = '  -0.17%  '
$c = $ws.Range("D20")
$c.Value = '''433.17'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.61%  '
$c = $ws.Range("D21")
$c.Value = '''13.45'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.97%  '
$c = $ws.Range("D22")
$c.Value = '''0.669'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.94%  '
$c = $ws.Range("D23")
$c.Value = '''7.03'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.54%  '
$c = $ws.Range("D24")
$c.Value = '''81.10'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$c = $ws.Range("D25")
$c.Value = '''10.67'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("E26").Value = '  -3.02%  '
$c = $ws.Range("D27")
$c.Value = '''11.71'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -7.73%  '
$c = $ws.Range("D31")
$c.Value = '''6.83'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.70%  '
$c = $ws.Range("D32")
$c.Value = '''26.17'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  -2.97%  '
$ws.Range("D35").Value = '0.0₃0857'
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("E36").Value = '  -2.26%  '
$c = $ws.Range("D37")
$c.Value = '''5.55'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.82%  '
$c = $ws.Range("D38")
$c.Value = '''49.20'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("E42").Value = '  -2.77%  '
$c = $ws.Range("D43")
$c.Value = '''0.273'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.06%  '
$c = $ws.Range("D44")
$c.Value = '''38.38'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -6.99%  '
$ws.Range("D45").Value = '2.681.46'
$ws.Range("E45").Value = '  -0.44%  '
$c = $ws.Range("D46")
$c.Value = '''133.50'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.27%  '
$c = $ws.Range("D47")
$c.Value = '''0.0335'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.94%  '
$c = $ws.Range("D48")
$c.Value = '''355.67'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.95%  '
$c = $ws.Range("D50")
$c.Value = '''22.66'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.20%  '
$ws.Range("E51").Value = '  -2.21%  '
